# "Adding moves up to L"
#
# Append a new data row (row 10) to Sheet1's Pokemon table, for Accelgor
# at level 100. The table already has a header row (row 1) and data rows
# 2-9; this extends it by one more row and grows the sheet dimension from
# A1:K9 to A1:K10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowNum = 10

# Plain text values - none of these look like numbers, so a direct
# assignment is stored as text without any extra formatting quirks.
$textValues = @{
    "A" = "Accelgor"
    "C" = "Adamant Orb"
    "D" = "Sticky Hold"
    "E" = "199 HP / 96 Atk / 136 Def / 79 SpA"
    "F" = "Serious"
    "H" = "Acid Spray"
    "I" = "Giga Drain"
    "J" = "Guard Split"
    "K" = "Feint"
}

foreach ($col in $textValues.Keys) {
    $ws.Range("$col$rowNum").Value = $textValues[$col]
}

# The "level" column (B) holds "100" as text (matching how the rest of
# this column is stored in the existing sheet), even though it looks
# numeric. Compute it with a TEXT() formula, then flatten the formula to
# its plain text result via Copy + PasteSpecial (values only) so the
# cell ends up as a plain text value instead of a live formula.
$b10 = $ws.Range("B10")
$b10.Formula = '=TEXT(100,"0")'
$b10.Copy()
$b10.PasteSpecial(-4163)
$excel.CutCopyMode = 0

# The "iv_spread" column (G) is present for this row but left blank.
$g10 = $ws.Range("G10")
$g10.Value = "'"
$g10.Style = "Normal"
